$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.991.87"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.955.53"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.61%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.98"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.07%  "
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4841"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2940"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.40%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07079"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +3.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.89"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "107.15"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.970.97"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07791"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.52%  "
$ws.Range("E14").Value = "  -0.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7012"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "279.45"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.67%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.998.87"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007827"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.28"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.214.07"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.40%  "
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.549"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.501"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.811"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.98"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.73"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.177"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1048"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.387"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -4.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.571"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.613"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.447"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04899"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.94%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7470"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -4.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.166"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.07%  "
$ws.Range("E37").Value = "  +0.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02000"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.687"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.538"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "77.74"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +8.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.120"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8984"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "109.23"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.90%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4445"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.970"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +5.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "991.70"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +5.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.279"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.49%  "
$ws.Range("E50").Value = "  -2.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "35.86"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.46%  "
